$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "correo" column (H) with header + generated addresses ---
$ws.Range("H1").Value = "correo"

$emails = @(
    "rsnaith0@symantec.com",
    "saucoate1@pcworld.com",
    "dcarefull2@wikimedia.org",
    "mtrathan3@ca.gov",
    "sluis4@vk.com",
    "saggas5@dagondesign.com",
    "apuncher6@linkedin.com",
    "dmarkovich7@trellian.com",
    "lescott8@timesonline.co.uk",
    "koharney9@businessinsider.com",
    "ssmailsa@livejournal.com",
    "sgatheralb@prweb.com",
    "glownesc@nasa.gov",
    "mgronousd@guardian.co.uk",
    "manscotte@sina.com.cn",
    "vtinkerf@people.com.cn",
    "cbertrandg@nymag.com",
    "knutkinh@google.com",
    "lkippeni@bing.com",
    "amcmillianj@constantcontact.com",
    "ehebbesk@google.com.hk",
    "pvardiel@stumbleupon.com",
    "hgobolosm@mozilla.org",
    "nmallenn@theatlantic.com",
    "tthwaiteo@irs.gov",
    "enursep@over-blog.com",
    "mferneq@webs.com"
)

for ($i = 0; $i -lt $emails.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $emails[$i]
}

# --- Data corrections ---
# Row 9 (cedula 8888, Jorge ...): clear the "segundo_nombre" (Manuel)
$ws.Range("C9").Value = ""

# Row 23 (cedula 1111, Mariita ...): fill in missing "direccion"
$ws.Range("F23").Value = "Calle 3 No. 2 - 3"

# Row 27 (cedula 8888): merge first+second name into primer_nombre, clear segundo_nombre
$ws.Range("B27").Value = "Jorge Manuel"
$ws.Range("C27").Value = ""

# --- Formatting: cedula column (A2:A28) stored as text ---
$ws.Range("A2:A28").NumberFormat = "@"

# Touch B28 so it carries its own explicit style (matches source edit)
$ws.Range("B28").Font.Name = $ws.Range("B28").Font.Name

# --- View state: scrolled down a bit, active cell on B28 ---
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("B28").Select()
